{"js": "// Update the worksheet date and all \"A\u00f7B=\" division expressions in the\n// table with new values (same document order as the source OOXML diff).\n//\n// Every source string in this document is unique, so exact-text search\n// is a safe way to locate each target run. However some NEW values\n// coincidentally equal some OLD values elsewhere in the document\n// (e.g. \"84\u00f72=\" becomes \"76\u00f72=\", while the ORIGINAL \"76\u00f72=\" cell\n// becomes \"74\u00f77=\", which in turn was itself an original value that\n// becomes \"36\u00f72=\"). A naive sequential search-then-replace loop would\n// therefore sometimes re-match a cell that a previous step just wrote,\n// corrupting the chain. To avoid that we resolve every search against\n// the pristine document first (phase 1) and only after all of them\n// have been located do we perform the text replacements (phase 2).\nconst pairs = [\n  [\"2024-10-03 Thursday\", \"2024-10-04 Friday\"],\n  [\"44\u00f74=\", \"63\u00f76=\"],\n  [\"19\u00f75=\", \"58\u00f79=\"],\n  [\"77\u00f78=\", \"79\u00f79=\"],\n  [\"28\u00f76=\", \"75\u00f79=\"],\n  [\"68\u00f79=\", \"70\u00f74=\"],\n  [\"68\u00f72=\", \"15\u00f75=\"],\n  [\"84\u00f75=\", \"20\u00f72=\"],\n  [\"62\u00f79=\", \"13\u00f79=\"],\n  [\"84\u00f72=\", \"76\u00f72=\"],\n  [\"75\u00f74=\", \"70\u00f72=\"],\n  [\"78\u00f78=\", \"17\u00f73=\"],\n  [\"24\u00f79=\", \"59\u00f77=\"],\n  [\"27\u00f73=\", \"30\u00f76=\"],\n  [\"48\u00f73=\", \"53\u00f75=\"],\n  [\"66\u00f75=\", \"18\u00f72=\"],\n  [\"22\u00f73=\", \"14\u00f79=\"],\n  [\"47\u00f78=\", \"67\u00f78=\"],\n  [\"73\u00f74=\", \"37\u00f74=\"],\n  [\"26\u00f75=\", \"63\u00f74=\"],\n  [\"77\u00f77=\", \"49\u00f77=\"],\n  [\"79\u00f75=\", \"61\u00f75=\"],\n  [\"59\u00f72=\", \"36\u00f78=\"],\n  [\"74\u00f77=\", \"36\u00f72=\"],\n  [\"65\u00f79=\", \"12\u00f74=\"],\n  [\"76\u00f72=\", \"74\u00f77=\"],\n];\n\n// Phase 1: locate every occurrence in the pristine document.\nconst searchResults = [];\nfor (const [oldText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  searchResults.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < pairs.length; i++) {\n  const found = searchResults[i].items.length;\n  if (found !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${pairs[i][0]}\", found ${found}`);\n  }\n}\n\n// Phase 2: replace each found range with its new text. Replacing text\n// inside the existing range (rather than rewriting the whole\n// paragraph/cell) preserves the surrounding run/paragraph formatting\n// (rFonts, sz, jc, etc.), matching the source diff exactly.\nfor (let i = 0; i < pairs.length; i++) {\n  const [, newText] = pairs[i];\n  searchResults[i].items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and all \"A\u00f7B=\" division expressions in the\n# table with new values.\n#\n# We address every target by POSITION (paragraph index / table cell\n# row+column) rather than by searching for its old text. Several of the\n# new values coincidentally equal OLD values that live elsewhere in the\n# document (e.g. \"84\u00f72=\" becomes \"76\u00f72=\", while the document's original\n# \"76\u00f72=\" cell becomes \"74\u00f77=\", which itself was an original value that\n# becomes \"36\u00f72=\"), so a naive Find/Replace sweep could re-match a cell\n# that an earlier step just wrote. Addressing cells by (row, column)\n# sidesteps that problem entirely, and assigning directly to\n# `Range.Text` replaces only the content (not the trailing paragraph /\n# cell-mark), which preserves the existing run/paragraph formatting\n# (rFonts, sz, jc, ...).\n\n$d = $word.ActiveDocument\n\n# First paragraph: the date line (outside of the table).\n$d.Paragraphs.Item(1).Range.Text = \"2024-10-04 Friday\"\n\n# The table has 20 rows x 5 columns; every 4th row (1, 5, 9, 13, 17 in\n# 1-based indexing) holds the five division expressions for that line,\n# the rows in between are blank spacer rows.\n$tbl = $d.Tables.Item(1)\n\n$newValues = @{\n  1  = @(\"63\u00f76=\", \"58\u00f79=\", \"79\u00f79=\", \"75\u00f79=\", \"70\u00f74=\");\n  5  = @(\"15\u00f75=\", \"20\u00f72=\", \"13\u00f79=\", \"76\u00f72=\", \"70\u00f72=\");\n  9  = @(\"17\u00f73=\", \"59\u00f77=\", \"30\u00f76=\", \"53\u00f75=\", \"18\u00f72=\");\n  13 = @(\"14\u00f79=\", \"67\u00f78=\", \"37\u00f74=\", \"63\u00f74=\", \"49\u00f77=\");\n  17 = @(\"61\u00f75=\", \"36\u00f78=\", \"36\u00f72=\", \"12\u00f74=\", \"74\u00f77=\");\n}\n\nforeach ($row in $newValues.Keys) {\n  $values = $newValues[$row]\n  for ($col = 1; $col -le $values.Length; $col++) {\n    $tbl.Cell($row, $col).Range.Text = $values[$col - 1]\n  }\n}\n"}
